$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.811.86'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').Value = '3.155.70'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '3.153.53'
$ws.Range('E8').Value = '  +0.53%  '
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.02'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.35%  '
$ws.Range('E12').Value = '  +1.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.69'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000249'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').Value = '3.675.75'
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('E17').Value = '  +3.34%  '
$ws.Range('D18').Value = '64.423.87'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').Value = '3.156.68'
$ws.Range('E19').Value = '  +0.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '476.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.10%  '
$ws.Range('E22').Value = '  +2.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.71'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '82.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.49%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.93'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.76%  '
$ws.Range('B29').Value = 'NEARProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.47'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.58%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.74'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.24'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.90%  '
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('E33').Value = '  +7.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.83'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.17%  '
$ws.Range('D35').Value = '0.0₃0879'
$ws.Range('E35').Value = '  +3.67%  '
$ws.Range('E36').Value = '  +6.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.08'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.24'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.63%  '
$ws.Range('E39').Value = '  +0.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '467.42'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.13%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '51.56'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.40'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.303'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.53%  '
$ws.Range('E44').Value = '  +3.00%  '
$ws.Range('D45').Value = '2.902.49'
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('E46').Value = '  +3.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '38.91'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.51'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.96'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.29'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.56%  '
